# Apply updated profit-calculation figures across all item sheets
# (values refreshed by the scheduled Sagittarius profits runner).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 6.470588
$ws.Range("I11").Value = 6.470588
$ws.Range("K11").Value = 6.470588
$ws.Range("M11").Value = 133.529412
$ws.Range("H17").Value = 6945.9443
$ws.Range("J17").Value = 7689.1875
$ws.Range("L17").Value = 23067.5625
$ws.Range("N17").Value = -23403.5625
$ws.Range("H40").Value = 1836.5
$ws.Range("J40").Value = 2059.6
$ws.Range("L40").Value = 2059.6
$ws.Range("N40").Value = -2409.6
$ws.Range("H97").Value = 1562.25
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1562.25
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 4686.75
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -5678.75
$ws.Range("H99").Value = 507.83334
$ws.Range("I99").Value = 274
$ws.Range("J99").Value = 624.75
$ws.Range("K99").Value = 822
$ws.Range("L99").Value = 1874.25
$ws.Range("M99").Value = 676
$ws.Range("N99").Value = -4870.25
$ws.Range("H101").Value = 20000380
$ws.Range("I101").Value = 33333500
$ws.Range("J101").Value = 699.5
$ws.Range("K101").Value = 100000500
$ws.Range("L101").Value = 2098.5
$ws.Range("M101").Value = -99998878
$ws.Range("N101").Value = -5342.5
$ws.Range("H138").Value = 2389.7778
$ws.Range("I138").Value = 2468.7273
$ws.Range("J138").Value = 2265.7144
$ws.Range("K138").Value = 7406.1819
$ws.Range("L138").Value = 6797.1432
$ws.Range("M138").Value = -2266.1819
$ws.Range("N138").Value = -17077.1432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 5555
$ws.Range("J21").Value = 5555
$ws.Range("L21").Value = 5555
$ws.Range("N21").Value = -6303
$ws.Range("H45").Value = 6900
$ws.Range("I45").Value = 6900
$ws.Range("K45").Value = 6900
$ws.Range("M45").Value = -6523
$ws.Range("H97").Value = 1003.9
$ws.Range("I97").Value = 1003.9
$ws.Range("K97").Value = 1003.9
$ws.Range("M97").Value = -507.9
$ws.Range("H102").Value = 866
$ws.Range("I102").Value = 881.3
$ws.Range("J102").Value = 789.5
$ws.Range("K102").Value = 881.3
$ws.Range("L102").Value = 789.5
$ws.Range("M102").Value = 740.7
$ws.Range("N102").Value = -4033.5
$ws.Range("H122").Value = 2804.5
$ws.Range("J122").Value = 2941.7144
$ws.Range("L122").Value = 8825.143199999999
$ws.Range("N122").Value = -13725.1432
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 66992
$ws.Range("J52").Value = 66992
$ws.Range("L52").Value = 66992
$ws.Range("N52").Value = -67518
$ws.Range("H105").Value = 4066.3333
$ws.Range("I105").Value = 3750
$ws.Range("K105").Value = 3750
$ws.Range("M105").Value = -2003
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H121").Value = 66992
$ws.Range("J121").Value = 66992
$ws.Range("L121").Value = 66992
$ws.Range("N121").Value = -70486

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1928.75
$ws.Range("I31").Value = 1681.6923
$ws.Range("J31").Value = 2999.3333
$ws.Range("K31").Value = 1681.6923
$ws.Range("L31").Value = 2999.3333
$ws.Range("M31").Value = -1386.6923
$ws.Range("N31").Value = -3589.3333
$ws.Range("H34").Value = 1928.75
$ws.Range("I34").Value = 1681.6923
$ws.Range("J34").Value = 2999.3333
$ws.Range("K34").Value = 1681.6923
$ws.Range("L34").Value = 2999.3333
$ws.Range("M34").Value = -1479.6923
$ws.Range("N34").Value = -3403.3333

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 69980.336
$ws.Range("J37").Value = 69980.336
$ws.Range("L37").Value = 209941.008
$ws.Range("N37").Value = -210165.008
$ws.Range("H52").Value = 1199
$ws.Range("J52").Value = 1199
$ws.Range("L52").Value = 3597
$ws.Range("N52").Value = -4129
$ws.Range("H75").Value = 8214.111000000001
$ws.Range("J75").Value = 10030.571
$ws.Range("L75").Value = 30091.713
$ws.Range("N75").Value = -32087.713
$ws.Range("H78").Value = 8214.111000000001
$ws.Range("J78").Value = 10030.571
$ws.Range("L78").Value = 90275.139
$ws.Range("N78").Value = -100259.139
$ws.Range("H87").Value = 2500
$ws.Range("I87").Value = 2500
$ws.Range("K87").Value = 7500
$ws.Range("M87").Value = -6252
$ws.Range("H90").Value = 2500
$ws.Range("I90").Value = 2500
$ws.Range("K90").Value = 22500
$ws.Range("M90").Value = -16260
$ws.Range("H132").Value = 2929.875
$ws.Range("I132").Value = 2012.5
$ws.Range("K132").Value = 18112.5
$ws.Range("M132").Value = -15582.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 650.4666999999999
$ws.Range("I97").Value = 740.5454999999999
$ws.Range("J97").Value = 402.75
$ws.Range("K97").Value = 740.5454999999999
$ws.Range("L97").Value = 402.75
$ws.Range("M97").Value = -244.5454999999999
$ws.Range("N97").Value = -1394.75
$ws.Range("H107").Value = 1990.7
$ws.Range("I107").Value = 1367.8
$ws.Range("K107").Value = 1367.8
$ws.Range("M107").Value = 552.2
$ws.Range("H132").Value = 1996
$ws.Range("I132").Value = 1996
$ws.Range("K132").Value = 5988
$ws.Range("M132").Value = -3458
$ws.Range("H136").Value = 8595
$ws.Range("J136").Value = 8595
$ws.Range("L136").Value = 25785
$ws.Range("N136").Value = -30885

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 229.55556
$ws.Range("I55").Value = 209.75
$ws.Range("J55").Value = 245.4
$ws.Range("K55").Value = 209.75
$ws.Range("L55").Value = 245.4
$ws.Range("M55").Value = -36.75
$ws.Range("N55").Value = -591.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 23857.143
$ws.Range("J62").Value = 21166.334
$ws.Range("L62").Value = 21166.334
$ws.Range("N62").Value = -22414.334
$ws.Range("H65").Value = 23857.143
$ws.Range("J65").Value = 21166.334
$ws.Range("L65").Value = 105831.67
$ws.Range("N65").Value = -112071.67
$ws.Range("H81").Value = 910506.4399999999
$ws.Range("I81").Value = 1556.9
$ws.Range("K81").Value = 3113.8
$ws.Range("M81").Value = -2052.8
$ws.Range("H84").Value = 910506.4399999999
$ws.Range("I84").Value = 1556.9
$ws.Range("K84").Value = 15569
$ws.Range("M84").Value = -10265
$ws.Range("H100").Value = 7143889
$ws.Range("I100").Value = 10001104
$ws.Range("K100").Value = 20002208
$ws.Range("M100").Value = -20001667
$ws.Range("H132").Value = 3843.1924
$ws.Range("I132").Value = 3580.25
$ws.Range("J132").Value = 6998.5
$ws.Range("K132").Value = 10740.75
$ws.Range("L132").Value = 20995.5
$ws.Range("M132").Value = -8210.75
$ws.Range("N132").Value = -26055.5
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
